$wb = $excel.ActiveWorkbook

# The "prereqs_test" sheet gets filled in with all prereq rows except the
# elective courses (CPSC-55200, CPSC-55500, CPSC-57100, CPSC-57200, CPSC-57400).
$wsTest = $wb.Worksheets.Item("prereqs_test")

$rows = @(
    @("CPSC-50100", "CPSC-51100"),
    @("CPSC-50100", "CPSC-51000"),
    @("CPSC-50100", "CPSC-53000"),
    @("CPSC-50100", "CPSC-54000"),
    @("MATH-51100", "MATH-51200"),
    @("MATH-51000", "CPSC-59000"),
    @("MATH-51100", "CPSC-59000"),
    @("MATH-51200", "CPSC-59000"),
    @("CPSC-50100", "CPSC-59000"),
    @("CPSC-51000", "CPSC-59000"),
    @("CPSC-51100", "CPSC-59000"),
    @("CPSC-53000", "CPSC-59000"),
    @("CPSC-54000", "CPSC-59000"),
    @("CPSC-55000", "CPSC-59000")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $wsTest.Cells.Item($r, 1).Value = $rows[$i][0]
    $wsTest.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Selection on "course_rotations" sheet changed from A2:A3 to A11:A18.
$wsCourses = $wb.Worksheets.Item("course_rotations")
$wsCourses.Range("A11:A18").Select() | Out-Null

# Selection on "prereqs" sheet changed (from topLeftCell scroll to a range selection).
$wsPrereqs = $wb.Worksheets.Item("prereqs")
$wsPrereqs.Range("A2:B23").Select() | Out-Null

# Leave the "prereqs_test" sheet active with its new selection, matching the
# saved workbook (prereqs_test is the tab that was visible/edited last).
$wsTest.Activate()
$wsTest.Range("E12").Select() | Out-Null
